$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simple price/volume updates: row => (D value, E value)
$updates = @{
    2  = @("64.398.87", "  +0.22%  ")
    3  = @("3.516.83", "  +0.48%  ")
    4  = @($null, "  +0.00%  ")
    5  = @("592.55", "  +1.44%  ")
    6  = @("134.72", "  -0.30%  ")
    7  = @($null, "  -0.02%  ")
    8  = @($null, "  +0.22%  ")
    9  = @("7.65", "  +7.25%  ")
    10 = @($null, "  +0.39%  ")
    11 = @($null, "  +3.91%  ")
    12 = @("4.116.00", "  +0.38%  ")
    13 = @($null, "  +1.34%  ")
    14 = @($null, "  +0.96%  ")
    15 = @("3.516.55", "  +0.30%  ")
    16 = @("25.87", "  +0.16%  ")
    17 = @("64.385.18", "  +0.17%  ")
    18 = @("10.01", "  +2.59%  ")
    19 = @($null, "  +3.51%  ")
    20 = @("13.59", "  -1.75%  ")
    21 = @("394.59", "  +2.93%  ")
    22 = @($null, "  +1.44%  ")
    23 = @("3.656.94", "  +0.46%  ")
    24 = @("74.72", "  +1.01%  ")
    25 = @($null, "  +0.16%  ")
    26 = @($null, "  +0.29%  ")
    27 = @("0.0000118", "  +3.23%  ")
    28 = @($null, "  +0.12%  ")
    29 = @("7.43", "  -1.13%  ")
    30 = @("2.27", "  +2.06%  ")
    31 = @("8.32", "  +0.53%  ")
    32 = @($null, "  -5.89%  ")
    33 = @($null, "  +6.98%  ")
    34 = @("3.547.71", "  +0.80%  ")
    36 = @("23.45", "  -0.33%  ")
    37 = @("5.34", "  +0.84%  ")
    38 = @("6.99", "  +2.37%  ")
    39 = @($null, "  +0.89%  ")
    40 = @("167.14", "  +1.80%  ")
    41 = @("0.0791", "  +1.11%  ")
    42 = @("0.813", "  +0.70%  ")
    45 = @("4.46", "  +1.32%  ")
    48 = @("6.84", "  +1.30%  ")
    49 = @("2.392.17", "  -3.23%  ")
    50 = @($null, "  -1.75%  ")
    51 = @($null, "  +0.74%  ")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($null -ne $vals[0]) {
        $ws.Cells.Item($row, 4).NumberFormat = "@"
        $ws.Cells.Item($row, 4).Value = $vals[0]
    }
    if ($null -ne $vals[1]) {
        $ws.Cells.Item($row, 5).NumberFormat = "@"
        $ws.Cells.Item($row, 5).Value = $vals[1]
    }
}

# Row 43 and 44 swap identities (EnergySwap <-> FirstDigitalUSD), with new D/E values
$ws.Cells.Item(43, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "1.00"
$ws.Cells.Item(43, 5).NumberFormat = "@"
$ws.Cells.Item(43, 5).Value = "  -0.03%  "

$ws.Cells.Item(44, 2).Value = "EnergySwap"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "25.39"
$ws.Cells.Item(44, 5).NumberFormat = "@"
$ws.Cells.Item(44, 5).Value = "  -3.21%  "

# Row 46 and 47 swap identities (Stacks <-> ONDO), with new D/E values
$ws.Cells.Item(46, 2).Value = "ONDO"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "1.18"
$ws.Cells.Item(46, 5).NumberFormat = "@"
$ws.Cells.Item(46, 5).Value = "  -1.67%  "

$ws.Cells.Item(47, 2).Value = "Stacks"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "1.66"
$ws.Cells.Item(47, 5).NumberFormat = "@"
$ws.Cells.Item(47, 5).Value = "  +2.54%  "
